# 2017-1-17: add some file
# Append two new subject columns (AC, AD) with a header row entry each
# plus their two data rows, matching the newly-added shared strings and
# sheet1 cell ranges in the commit's diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (become new shared-string entries 28 and 29)
$ws.Range("AC1").Value = "wnb-调节6Hz_20161230_113123_ASIC_EEG"
$ws.Range("AD1").Value = "wnb-调节6Hz_20170110_113300_ASIC_EEG"

# New probability values for the two new subjects/columns
$ws.Range("AC2").Value = 0.7138263665594855
$ws.Range("AD2").Value = 0.78640776699029125

$ws.Range("AC3").Value = 0.70186335403726707
$ws.Range("AD3").Value = 0.82935153583617749

# Keep the sheet's selection in sync with the now-wider used range
$null = $ws.Range("A1:AD3").Select()
